$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values (order matters: drives sharedStrings index assignment) ---
$ws.Range("B2").Value = "STT"
$ws.Range("C2").Value = "TÊN"

$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "use case"

$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "đặc tả"

$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "stackholder"

$ws.Range("B6").Value = 4
$ws.Range("C6").Value = "yêu cầu kỹ thuật"

$ws.Range("B7").Value = 5
$ws.Range("C7").Value = "thuật ngữ trong phần mềm"

$ws.Range("B8").Value = 6
$ws.Range("C8").Value = "testcase"

$ws.Range("B9").Value = 7
$ws.Range("C9").Value = "test plan"

$ws.Range("B10").Value = 8
$ws.Range("C10").Value = "đánh giá rủi ro"

$ws.Range("B11").Value = 9
$ws.Range("C11").Value = "quy trình nghiệp vụ"

$ws.Range("B12").Value = 10
$ws.Range("C12").Value = "giới thiệu phần mềm"

$ws.Range("B13").Value = 11
$ws.Range("C13").Value = "prototype"

$ws.Range("B14").Value = 12
$ws.Range("C14").Value = "các phần lấy yêu cầu khách hàng"

$ws.Range("B15").Value = 13
$ws.Range("C15").Value = "quản lý kế hoạch dự án"

# --- Column width ---
$ws.Columns.Item(3).ColumnWidth = 25.3

# --- Row heights for the thick top/bottom border rows ---
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(15).RowHeight = 15.75
$ws.Rows.Item(16).RowHeight = 15.75

# --- Header row fill (yellow) ---
$headerRange = $ws.Range("B2:C2")
$headerRange.Interior.Color = 65535

# --- Number column (B) fill (theme accent, light) for data rows ---
$numRange = $ws.Range("B3:B15")
$numRange.Interior.ThemeColor = 3
$numRange.Interior.TintAndShade = 0.79998168889431442

# --- Borders: outer double border around B2:C15, thin inner gridlines ---
$tableRange = $ws.Range("B2:C15")

# Outer border = double
$tableRange.Borders.Item(7).LineStyle = -4119   # xlEdgeLeft
$tableRange.Borders.Item(8).LineStyle = -4119   # xlEdgeTop
$tableRange.Borders.Item(9).LineStyle = -4119   # xlEdgeBottom
$tableRange.Borders.Item(10).LineStyle = -4119  # xlEdgeRight

# Inner gridlines = thin
$tableRange.Borders.Item(11).LineStyle = 1      # xlInsideVertical
$tableRange.Borders.Item(12).LineStyle = 1      # xlInsideHorizontal

# --- Page setup ---
$ws.PageSetup.Orientation = 1

# --- Selection (matches authored view state) ---
[void]$ws.Range("G17").Select()
